$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two trailing cells of row 7 (PriceChange, UpDown)
$ws.Range("X7").Value = 0.22000199999999381
$ws.Range("Y7").Value = "Up"

# Append a brand-new data row 8
$ws.Range("A8").Value = 42649.87908564815
$ws.Range("A8").NumberFormat = "m/d/yy h:mm"
$ws.Range("B8").Value = -2
$ws.Range("C8").Value = "Neutral"
$ws.Range("D8").Value = 34
$ws.Range("E8").Value = 13261
$ws.Range("F8").Value = 2279
$ws.Range("G8").Value = 56
$ws.Range("H8").Value = 40
$ws.Range("I8").Value = 87
$ws.Range("J8").Value = 12
$ws.Range("K8").Value = 15052
$ws.Range("L8").Value = 258
$ws.Range("M8").Value = 184
$ws.Range("N8").Value = 91
$ws.Range("O8").Value = 13
$ws.Range("P8").Value = "Noun"
$ws.Range("Q8").Value = 47.418521827693588
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = -0.089899999999999994
$ws.Range("S8").NumberFormat = "0.00%"
$ws.Range("T8").Value = -0.021000000000000001
$ws.Range("T8").NumberFormat = "0.00%"
$ws.Range("U8").Value = 6.65
$ws.Range("V8").Value = 1.88
$ws.Range("W8").Value = -2
